$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
# Fix the "Market Capitalization" header (now "Market Captilization")
$ws.Range("C1").Value = "Market Captilization"

# The bold/centered header style is no longer used - make B1:D1 use the
# same formatting as the data columns below them (copy number format /
# font / fill / border from row 2, which is what the columns already use).
$ws.Range("B2").Copy()
$ws.Range("B1").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C1").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("D1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Data rows --------------------------------------------------------
$ws.Range("B2").Value = 126.05
$ws.Range("C2").Value = 37232144384
$ws.Range("D2").Value = 158

$ws.Range("B3").Value = 16.79
$ws.Range("C3").Value = 10969948160
$ws.Range("D3").Value = 1191

$ws.Range("B4").Value = 73.17
$ws.Range("C4").Value = 4349510144
$ws.Range("D4").Value = 273

$ws.Range("B5").Value = 195.83
$ws.Range("C5").Value = 3080151367680
$ws.Range("D5").Value = 102

$ws.Range("B6").Value = 150.85
$ws.Range("C6").Value = 266143170560
$ws.Range("D6").Value = 132
